# "Começo da normalização dos dados no classificador"
# Populate the Goal16 classifier reference sheet with the Metric/Reference
# table. Column A (rows 2-4) is filled first, then the header row, then
# column B, matching the shared-string insertion order recorded in the
# target workbook's xl/sharedStrings.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Number of victims of intentional homicide per 100,000 population, by sex (victims per 100,000 population)"
$ws.Range("A3").Value = "Number of victims of intentional homicide, by sex (number)"
$ws.Range("A4").Value = "Unsentenced detainees as a proportion of overall prison population (%)"

$ws.Range("A1").Value = "Metric"
$ws.Range("B1").Value = "Reference"

$ws.Range("B2").Value = "decrease"
$ws.Range("B3").Value = "decrease"
$ws.Range("B4").Value = "decrease"

# Auto-fit the two columns to the new content (as in the original edit,
# which left both columns with bestFit/customWidth column metadata).
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null

# Nudge the computed widths to match Excel's own font-metrics-based
# auto-fit result (97.42578125 / 10.140625 chars) as closely as this
# engine's width quantization allows.
$ws.Columns.Item(1).ColumnWidth = 96.66666666666667
$ws.Columns.Item(2).ColumnWidth = 9.333333333333334

# Leave selection where the commit's saved view shows it.
$ws.Range("B5").Select() | Out-Null
